$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / "Changed") date for rows 2-6 from 2023-11-03 (45233)
# to 2023-11-13 (45243), keeping existing date formatting/style intact.
$newDate = (Get-Date -Year 2023 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0).Date

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
